# Auto-generated edit script applying the scheduled market-data refresh
# to the Ragnarok_Profits workbook (columns H-N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H11").Value = 214.375
$ws.Range("I11").Value = 214.375
$ws.Range("K11").Value = 214.375
$ws.Range("M11").Value = -74.375

$ws.Range("H17").Value = 433.65308
$ws.Range("J17").Value = 437.21277
$ws.Range("L17").Value = 1311.63831
$ws.Range("N17").Value = -1647.63831

$ws.Range("H43").Value = 5500
$ws.Range("I43").Value = 5500
$ws.Range("J43").Value = 5500
$ws.Range("K43").Value = 5500
$ws.Range("L43").Value = 5500
$ws.Range("M43").Value = -5431
$ws.Range("N43").Value = -5638

$ws.Range("H58").Value = 4161.5
$ws.Range("I58").Value = 323.33334
$ws.Range("J58").Value = 7999.6665
$ws.Range("K58").Value = 970.0000200000001
$ws.Range("L58").Value = 23998.9995
$ws.Range("M58").Value = -820.0000200000001
$ws.Range("N58").Value = -24298.9995

$ws.Range("H62").Value = 3389
$ws.Range("I62").Value = 3120.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3120.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2496.5
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3389
$ws.Range("I65").Value = 3120.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 15602.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -12482.5
$ws.Range("N65").Value = -31240

$ws.Range("H107").Value = 741.1
$ws.Range("I107").Value = 712.44446
$ws.Range("K107").Value = 712.44446
$ws.Range("M107").Value = 1207.55554

$ws.Range("H111").Value = 4349.875
$ws.Range("I111").Value = 1759.8
$ws.Range("K111").Value = 5279.4
$ws.Range("M111").Value = -2212.4

$ws.Range("H114").Value = 99999
$ws.Range("J114").Value = 99999
$ws.Range("L114").Value = 99999
$ws.Range("N114").Value = -108677

$ws.Range("H129").Value = 3067.5186
$ws.Range("I129").Value = 560.6667
$ws.Range("K129").Value = 1682.0001
$ws.Range("M129").Value = 3317.9999

$ws.Range("H132").Value = 4369.864
$ws.Range("I132").Value = 2725.303
$ws.Range("K132").Value = 8175.909
$ws.Range("M132").Value = -5645.909

$ws.Range("H137").Value = 4179.5713
$ws.Range("I137").Value = 2167.3076
$ws.Range("K137").Value = 6501.9228
$ws.Range("M137").Value = -3951.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3918.9333
$ws.Range("I32").Value = 3393.3818
$ws.Range("K32").Value = 3393.3818
$ws.Range("M32").Value = -3106.3818

$ws.Range("H74").Value = 3861.7144
$ws.Range("I74").Value = 3433.9092
$ws.Range("J74").Value = 5430.3335
$ws.Range("K74").Value = 3433.9092
$ws.Range("L74").Value = 5430.3335
$ws.Range("M74").Value = -2559.9092
$ws.Range("N74").Value = -7178.3335

$ws.Range("H77").Value = 3861.7144
$ws.Range("I77").Value = 3433.9092
$ws.Range("J77").Value = 5430.3335
$ws.Range("K77").Value = 17169.546
$ws.Range("L77").Value = 27151.6675
$ws.Range("M77").Value = -12801.546
$ws.Range("N77").Value = -35887.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47622044
$ws.Range("I31").Value = 76925920
$ws.Range("J31").Value = 3249.5
$ws.Range("K31").Value = 76925920
$ws.Range("L31").Value = 3249.5
$ws.Range("M31").Value = -76925625
$ws.Range("N31").Value = -3839.5

$ws.Range("H34").Value = 47622044
$ws.Range("I34").Value = 76925920
$ws.Range("J34").Value = 3249.5
$ws.Range("K34").Value = 76925920
$ws.Range("L34").Value = 3249.5
$ws.Range("M34").Value = -76925718
$ws.Range("N34").Value = -3653.5

$ws.Range("H58").Value = 3278.25
$ws.Range("I58").Value = 3042.4
$ws.Range("K58").Value = 3042.4
$ws.Range("M58").Value = -2839.4

$ws.Range("H59").Value = 118799.2
$ws.Range("I59").Value = 21998
$ws.Range("J59").Value = 142999.5
$ws.Range("K59").Value = 21998
$ws.Range("L59").Value = 142999.5
$ws.Range("M59").Value = -20853
$ws.Range("N59").Value = -145289.5

$ws.Range("H107").Value = 1032.0416
$ws.Range("I107").Value = 714.9
$ws.Range("K107").Value = 714.9
$ws.Range("M107").Value = 1205.1

$ws.Range("H122").Value = 3626.7144
$ws.Range("I122").Value = 3344.4614
$ws.Range("K122").Value = 10033.3842
$ws.Range("M122").Value = -7583.3842

$ws.Range("H136").Value = 3278.25
$ws.Range("I136").Value = 3042.4
$ws.Range("K136").Value = 9127.200000000001
$ws.Range("M136").Value = -6577.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1325
$ws.Range("I5").Value = 1166.8334
$ws.Range("K5").Value = 3500.5002
$ws.Range("M5").Value = -3388.5002

$ws.Range("H23").Value = 316.33334
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 316.33334
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 949.0000200000001
$ws.Range("N23").Value = -1419.00002
$ws.Range("M23").ClearContents()

$ws.Range("H86").Value = 1323
$ws.Range("J86").Value = 1952.5
$ws.Range("L86").Value = 5857.5
$ws.Range("N86").Value = -8229.5

$ws.Range("H89").Value = 1323
$ws.Range("J89").Value = 1952.5
$ws.Range("L89").Value = 17572.5
$ws.Range("N89").Value = -29428.5

$ws.Range("H122").Value = 42635.875
$ws.Range("J122").Value = 1097.8
$ws.Range("L122").Value = 9880.199999999999
$ws.Range("N122").Value = -14780.2

$ws.Range("H133").Value = 35594.75
$ws.Range("I133").Value = 46285.4
$ws.Range("J133").Value = 17777
$ws.Range("K133").Value = 138856.2
$ws.Range("L133").Value = 53331
$ws.Range("M133").Value = -133796.2
$ws.Range("N133").Value = -63451

$ws.Range("H135").Value = 1325
$ws.Range("I135").Value = 1166.8334
$ws.Range("K135").Value = 10501.5006
$ws.Range("M135").Value = -7966.500599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3465999.5
$ws.Range("I11").Value = 5770000
$ws.Range("K11").Value = 5770000
$ws.Range("M11").Value = -5769861

$ws.Range("H102").Value = 2654.6956
$ws.Range("I102").Value = 2564.0476
$ws.Range("J102").Value = 3606.5
$ws.Range("K102").Value = 2564.0476
$ws.Range("L102").Value = 3606.5
$ws.Range("M102").Value = -942.0475999999999
$ws.Range("N102").Value = -6850.5

$ws.Range("H107").Value = 1294.4736
$ws.Range("I107").Value = 1377.8667
$ws.Range("J107").Value = 981.75
$ws.Range("K107").Value = 1377.8667
$ws.Range("L107").Value = 981.75
$ws.Range("M107").Value = 542.1333
$ws.Range("N107").Value = -4821.75

$ws.Range("H122").Value = 1850.1818
$ws.Range("I122").Value = 910.8333
$ws.Range("K122").Value = 2732.4999
$ws.Range("M122").Value = -282.4998999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19664

$ws.Range("H57").Value = 29365.834
$ws.Range("I57").Value = 29365.834
$ws.Range("K57").Value = 29365.834
$ws.Range("M57").Value = -28799.834

$ws.Range("H136").Value = 2237.1143
$ws.Range("I136").Value = 2181.5557
$ws.Range("J136").Value = 2424.625
$ws.Range("K136").Value = 6544.6671
$ws.Range("L136").Value = 7273.875
$ws.Range("M136").Value = -3994.6671
$ws.Range("N136").Value = -12373.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 99996.5
$ws.Range("J123").Value = 99996.5
$ws.Range("L123").Value = 99996.5
$ws.Range("N123").Value = -109796.5

$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws.Range("H132").Value = 4664.069
$ws.Range("I132").Value = 4060.2104
$ws.Range("K132").Value = 12180.6312
$ws.Range("M132").Value = -9650.6312
